# Apply updated cryptocurrency Price (column D) and Volume(1h) (column E)
# figures, refreshed by the scheduled GitHub Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.447.63'
$ws.Range("E2").Value = '  -0.65%  '
$ws.Range("D3").Value = '1.823.46'
$ws.Range("E3").Value = '  -2.29%  '
$ws.Range("D4").Value = '''1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = '''332.66'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.73%  '
$ws.Range("D6").Value = '''1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.11%  '
$ws.Range("D7").Value = '''0.4585'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.66%  '
$ws.Range("D8").Value = '''0.3813'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.83%  '
$ws.Range("D9").Value = '''45.99'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.08%  '
$ws.Range("D10").Value = '''0.07857'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.45%  '
$ws.Range("D11").Value = '''0.9587'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.28%  '
$ws.Range("D12").Value = '''21.03'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.20%  '
$ws.Range("D13").Value = '1.831.45'
$ws.Range("E13").Value = '  -2.19%  '
$ws.Range("D14").Value = '''5.838'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.09%  '
$ws.Range("D15").Value = '''7.083'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.90%  '
$ws.Range("E16").Value = '  -1.16%  '
$ws.Range("D17").Value = '''89.41'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.89%  '
$ws.Range("D18").Value = '''0.06583'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.23%  '
$ws.Range("D19").Value = '''0.00001020'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.27%  '
$ws.Range("D20").Value = '''17.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.70%  '
$ws.Range("D21").Value = '''1.002'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.16%  '
$ws.Range("D22").Value = '27.436.85'
$ws.Range("E22").Value = '  -0.75%  '
$ws.Range("D23").Value = '''5.291'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = '''10.83'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.87%  '
$ws.Range("D25").Value = '''2.258'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.30%  '
$ws.Range("D26").Value = '2.045.33'
$ws.Range("E26").Value = '  -2.35%  '
$ws.Range("D27").Value = '''156.91'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.52%  '
$ws.Range("D28").Value = '''19.34'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.65%  '
$ws.Range("D29").Value = '''2.043'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.29%  '
$ws.Range("D30").Value = '''5.266'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.18%  '
$ws.Range("D31").Value = '''117.78'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.26%  '
$ws.Range("D32").Value = '''0.09306'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.58%  '
$ws.Range("D33").Value = '''0.9315'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.59%  '
$ws.Range("D34").Value = '''3.562'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.82%  '
$ws.Range("D35").Value = '''5.219'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.60%  '
$ws.Range("D36").Value = '''1.316'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.78%  '
$ws.Range("D37").Value = '''0.05915'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.90%  '
$ws.Range("D38").Value = '''0.02187'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.17%  '
$ws.Range("D39").Value = '''8.099'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.41%  '
$ws.Range("D40").Value = '''1.002'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.08%  '
$ws.Range("D41").Value = '''1.140'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.77%  '
$ws.Range("D42").Value = '''0.5751'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.33%  '
$ws.Range("D43").Value = '''0.1816'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.09%  '
$ws.Range("D44").Value = '''9.931'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = '''1.275'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.82%  '
$ws.Range("D46").Value = '''0.5394'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.16%  '
$ws.Range("D47").Value = '''11.81'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.09%  '
$ws.Range("D48").Value = '''1.873'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.71%  '
$ws.Range("D49").Value = '''110.35'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.28%  '
$ws.Range("D50").Value = '''0.06572'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.58%  '
$ws.Range("D51").Value = '''1.003'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -33.40%  '
